$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.09
$ws.Range("C2").Value = 13.72
$ws.Range("D2").Value = 10.7

$ws.Range("B3").Value = 13.87
$ws.Range("C3").Value = 13.57
$ws.Range("D3").Value = 10.59
